$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the 9 columns that change (B, C, D, E, G, I, K, L, O)
$cols = @("B","C","D","E","G","I","K","L","O")

# New values per row (row 2 .. row 25), in the order of $cols above
$data = @{
    2 = @(12.51331208847585, 6.912263319083058, 6.010673508816727, 11.2908506909073, 3.677292287183406, 27.07832644793021, 10.07852797806205, 9.922229151393081, 27.76632933286209)
    3 = @(12.27345099989719, 6.85066907514227, 5.896362037508655, 11.30014276568276, 3.679400724005145, 27.15935734365117, 9.912143914538737, 9.909755547017784, 27.8309406438944)
    4 = @(12.12626420580782, 6.811950908263359, 5.82681724138734, 11.30780336865601, 3.680764141029002, 27.2140637953306, 9.810410113506169, 9.903850147134868, 27.87623621159289)
    5 = @(12.06639090309262, 6.795952744753104, 5.798681945716044, 11.31141702323004, 3.681337107284234, 27.23760076930635, 9.769113475634041, 9.901886704145761, 27.89610586979976)
    6 = @(12.05645790708305, 6.793283135892049, 5.794023733116425, 11.31204678138432, 3.681433298189331, 27.24158412384024, 9.762267474662064, 9.901587492320788, 27.89949036592662)
    7 = @(12.12545618999523, 6.811736035175131, 5.82643691269268, 11.30785011177926, 3.680771797877935, 27.21437619022705, 9.809852451059808, 9.903821870925338, 27.87649847021361)
    8 = @(12.43063953578097, 6.891215716808732, 5.971151240928883, 11.29364897150953, 3.678005024238602, 27.10523677355583, 10.02110210488156, 9.917565689362135, 27.78743841102485)
    9 = @(13.02620325388184, 7.039690231349217, 6.258166287298322, 11.28129950215859, 3.67312298691316, 26.93059715896635, 10.43643056774162, 9.958325000708806, 27.65753530248218)
    10 = @(13.45734344318468, 7.143946235302605, 6.468563476636622, 11.28164801966347, 3.669864011970881, 26.82640665298305, 10.73919152387102, 9.996531418342224, 27.58952486286319)
    11 = @(13.65116375375089, 7.19025976018058, 6.563685184763663, 11.28384362880539, 3.668451859399855, 26.7842665331104, 10.87579471230153, 10.01566696411986, 27.56456895193388)
    12 = @(13.72415480047029, 7.207631901482809, 6.599582058443505, 11.2849669382473, 3.667927176298742, 26.76906678640696, 10.92731245541366, 10.02316170505496, 27.55598086473502)
    13 = @(13.70845397605844, 7.203897963023707, 6.591857122640306, 11.28471205045325, 3.668039729202607, 26.77230659103677, 10.91622732898133, 10.02153658991785, 27.55779209009617)
    14 = @(13.65717727947773, 7.191692331577398, 6.56664113629804, 11.28393020102485, 3.668408491909449, 26.78300084274678, 10.88003760941558, 10.01627860333447, 27.56384511562423)
    15 = @(13.62571405787007, 7.184194272213107, 6.551178391588104, 11.28348927316156, 3.668635679417932, 26.78965012381461, 10.85784146920583, 10.01309018028413, 27.56766510104839)
    16 = @(13.44462409053901, 7.140896646904743, 6.462331821513922, 11.2815454366053, 3.669957711156691, 26.82926655086017, 10.7302373011104, 9.995315860333392, 27.5912763414351)
    17 = @(13.33288724605749, 7.114045753479531, 6.407648125240321, 11.28087400396029, 3.670786722286215, 26.8549176829118, 10.65163259548866, 9.984858781156493, 27.60729475764253)
    18 = @(13.26840606766382, 7.098497397547875, 6.376142276706967, 11.28067962509579, 3.671270174276335, 26.87016620932582, 10.60631823859732, 9.9790096083813, 27.61707112008392)
    19 = @(13.24653952971994, 7.093215232018048, 6.365467014203376, 11.28064678479234, 3.671435002514517, 26.87541399976844, 10.5909594200972, 9.977057711030415, 27.62047785999921)
    20 = @(13.34480441851544, 7.116914927240016, 6.413475095949826, 11.28092563442271, 3.670697787166454, 26.85213586517338, 10.6600112253309, 9.985954857560793, 27.60553128739419)
    21 = @(13.67225002496236, 7.195281965696074, 6.574051339066951, 11.28415193678078, 3.668299904512582, 26.77983910054487, 10.89067351313733, 10.01781628880554, 27.56204378229248)
    22 = @(13.88386616238892, 7.245530480096554, 6.678259834861668, 11.28796120946788, 3.666791412760139, 26.73700685497845, 11.04017399849958, 10.04008618155674, 27.53864828002277)
    23 = @(13.77116413246923, 7.218802408498176, 6.622721609787948, 11.28577290353239, 3.667591172229045, 26.75946239032035, 10.96051290288124, 10.02806930710425, 27.55067445407883)
    24 = @(13.33941741360189, 7.115618120418196, 6.41084092965815, 11.28090169529244, 3.670737973416718, 26.85339196340655, 10.65622362756709, 9.9854588142097, 27.60632678572559)
    25 = @(12.86588994426389, 7.000348643507673, 6.180432145055921, 11.28298325937814, 3.67438587827094, 26.97361414723076, 10.32427710478484, 9.945835383544514, 27.68787070493528)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
